$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Acetate+CO2+H2" ---
# Flip the sign of the two "pulled out" ATP / Fd reaction flux values, and
# leave the selection sitting on those two cells (this sheet is no longer
# the active tab after the edit).
$ws1 = $wb.Worksheets.Item("Acetate+CO2+H2")
$ws1.Range("C21").Value = 961.80277100000001
$ws1.Range("C22").Value = -1000
$ws1.Activate() | Out-Null
$ws1.Range("C21:C22").Select() | Out-Null

# --- Sheet 2: "Acetate+CO2+Formate" ---
# Same value flips, plus this sheet becomes the active tab, scrolled so
# row 6 is at the top, with D32:D35 selected (D32 active).
$ws2 = $wb.Worksheets.Item("Acetate+CO2+Formate")
$ws2.Range("C21").Value = 961.80277100000001
$ws2.Range("C22").Value = -1000
$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$ws2.Range("D32:D35").Select() | Out-Null
